# Rename sheets and round displayed data values as part of creating the
# "compare tab" beginning.

$wb = $excel.ActiveWorkbook

# Rename the worksheets (strip the redundant "mat_mul_" prefix and
# normalize the ccm/ram casing to CCM/RAM).
$wb.Worksheets.Item(1).Name = "data_CCM code_FLASH"
$wb.Worksheets.Item(2).Name = "data_CCM code_CCM"
$wb.Worksheets.Item(3).Name = "data_RAM code_FLASH"
$wb.Worksheets.Item(4).Name = "data_RAM code_CCM"

# Round the "intensity" row (row 2) to whole numbers and convert the
# "energy" row (row 5) from joules to kilojoules, rounded to 3 decimals.
foreach ($ws in $wb.Worksheets) {
    foreach ($col in @("B", "C", "D")) {
        $intensityCell = $ws.Range($col + "2")
        $intensityCell.Value2 = [math]::Round([double]$intensityCell.Value2, 0)

        $energyCell = $ws.Range($col + "5")
        $energyCell.Value2 = [math]::Round([double]$energyCell.Value2 / 1000.0, 3)
    }
}
